# Daily attendance processing - 2026-01-13 22:36:34
# Reorders the "Recorded By" (column G) entries for each session row so
# that any "System" entry is moved to the front of the comma separated
# list (preserving the relative order of the remaining entries). When a
# cell has no "System" entry but still lists multiple recorders, the
# order of the entries is reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notmatch ",") { continue }

    $parts = @($val -split "," | ForEach-Object { $_.Trim() })

    # Case-sensitive check for an exact "System" token (PowerShell's default
    # -eq/-ne/-contains operators are case-insensitive, which would wrongly
    # also match the distinct lowercase "system" entries).
    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) { $rest += $p }
        }
        $newParts = @("System") + $rest
    } else {
        $newParts = $parts[($parts.Count - 1)..0]
    }

    $newVal = ($newParts -join ", ")

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
